{"js": "// Word templating edit: turn the single \"Title1 / Title2\" block (which used\n// AppendParagraph() to fake a second row) into a real nested\n// `foreach(var test in project.StringList)` loop that prints `test` once per\n// iteration, replacing the old \"Title2: <%= project.ProjectName %>\" line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two target paragraphs by their known (pre-edit) text so the\n// script is resilient to exact index assumptions.\nconst items = paragraphs.items;\nlet title1Para = null;\nlet title2Para = null;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (title1Para === null && t.indexOf(\"Title1:\") !== -1) {\n    title1Para = items[i];\n  } else if (title2Para === null && t.indexOf(\"Title2:\") !== -1) {\n    title2Para = items[i];\n  }\n}\n\nif (!title1Para || !title2Para) {\n  throw new Error(\"Could not locate Title1/Title2 template paragraphs\");\n}\n\n// --- Paragraph \"Title1: <%= project.ProjectName %><% AppendParagraph(); %>\"\n// Drop the trailing \"<% AppendParagraph(); %>\" run so the paragraph just\n// reads \"Title1: <%= project.ProjectName %>\".\nlet found = title1Para.search(\"<% AppendParagraph(); %>\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].delete();\n  await context.sync();\n}\n\n// Insert the new \"foreach\" opening line right after the Title1 paragraph.\ntitle1Para.insertParagraph(\n  \"<% foreach(var test in project.StringList) { %>\",\n  \"After\"\n);\nawait context.sync();\n\n// --- Paragraph \"Title2:\" + bookmark + \" <%= project.ProjectName %><% AppendParagraph(); %>\"\n// Strip everything but the bookmark, then write \"<%= test %>\" ahead of it so\n// the paragraph becomes \"<%= test %>\" followed by the (untouched) bookmark.\nfound = title2Para.search(\"Title2:\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].delete();\n  await context.sync();\n}\n\nfound = title2Para.search(\" <%= project.ProjectName %>\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].delete();\n  await context.sync();\n}\n\nfound = title2Para.search(\"<% AppendParagraph(); %>\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].delete();\n  await context.sync();\n}\n\nconst startRange = title2Para.getRange(\"Start\");\nstartRange.insertText(\"<%= test %>\", \"Before\");\nawait context.sync();\n\n// Insert the loop-closing \"<% } %>\" line right after the (now bookmark-only\n// plus \"<%= test %>\") paragraph.\ntitle2Para.insertParagraph(\"<% } %>\", \"After\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction ReplaceInRange($range, $searchText, $replaceText) {\n    $f = $range.Find\n    $f.Text = $searchText\n    $f.Replacement.Text = $replaceText\n    $f.Execute($f.Text, $false, $false, $false, $false, $false, $true, 1, $false, $f.Replacement.Text, 2) | Out-Null\n}\n\n# Locate the two template paragraphs (\"Title1: ...\" and \"Title2: ...\") by\n# their current text instead of hard-coded indexes, so the script is\n# resilient to the exact paragraph numbering.\n$title1Index = -1\n$title2Index = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($title1Index -eq -1 -and $t -like \"*Title1:*\") {\n        $title1Index = $i\n    } elseif ($title2Index -eq -1 -and $t -like \"*Title2:*\") {\n        $title2Index = $i\n    }\n}\n\n# --- \"Title1: <%= project.ProjectName %><% AppendParagraph(); %>\"\n# Drop the trailing \"<% AppendParagraph(); %>\" run.\n$p1Range = $d.Paragraphs.Item($title1Index).Range\nReplaceInRange $p1Range \"<% AppendParagraph(); %>\" \"\"\n\n# Insert the new \"foreach\" opening line right after the Title1 paragraph.\n$d.Paragraphs.Item($title1Index).Range.InsertParagraphAfter()\n$title2Index = $title2Index + 1\n$d.Paragraphs.Item($title1Index + 1).Range.Text = \"<% foreach(var test in project.StringList) { %>\"\n\n# --- \"Title2:\" + bookmark + \" <%= project.ProjectName %><% AppendParagraph(); %>\"\n# Strip everything but the bookmark.\n$p2Range = $d.Paragraphs.Item($title2Index).Range\nReplaceInRange $p2Range \"Title2:\" \"\"\nReplaceInRange $p2Range \" <%= project.ProjectName %>\" \"\"\nReplaceInRange $p2Range \"<% AppendParagraph(); %>\" \"\"\n\n# Write \"<%= test %>\" ahead of the (now bookmark-only) paragraph content.\n$d.Paragraphs.Item($title2Index).Range.InsertBefore(\"<%= test %>\")\n\n# Insert the loop-closing \"<% } %>\" line right after that paragraph.\n$d.Paragraphs.Item($title2Index).Range.InsertParagraphAfter()\n$d.Paragraphs.Item($title2Index + 1).Range.Text = \"<% } %>\"\n"}
